$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 239
$ws.Range("I2").Value = 605
$ws.Range("J2").Value = 2493
$ws.Range("L2").Value = 678
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 438
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 8
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 32
$ws.Range("S2").Value = 288
$ws.Range("T2").Value = 420
$ws.Range("U2").Value = 32
$ws.Range("V2").Value = 3651
$ws.Range("X2").Value = 3896
$ws.Range("Z2").Value = 52
$ws.Range("AA2").Value = 24
